# Apply updated values to Sheet1 per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 13.377
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.694
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("A18").Value = -21.694
$ws.Range("E20").Value = 12.932
